# Updated cryptos list on Wed Dec 13 04:49:05 UTC 2023 with GitHub Actions
# Refreshes Price (column D) / Volume(1h) (column E) figures for each coin row,
# and fixes the ranking swap between Avalanche/Dogecoin (rows 11-12) and
# BinanceUSD/Cronos (rows 48-49). Price values that look like plain numbers are
# prefixed with a leading apostrophe so Excel keeps storing them as text
# (matching the original inline-string cells) instead of converting them to
# numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.854.57"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "2.165.73"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'246.71"
$ws.Range("E5").Value = "  -2.65%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").Value = "'66.40"
$ws.Range("E7").Value = "  -6.11%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.568"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").Value = "'58.29"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0925"
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "'35.60"
$ws.Range("E12").Value = "  -14.73%  "
$ws.Range("D13").Value = "'0.104"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "'6.89"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "2.484.37"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "'0.861"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'14.25"
$ws.Range("E17").Value = "  -4.28%  "
$ws.Range("D18").Value = "2.150.38"
$ws.Range("E18").Value = "  -3.59%  "
$ws.Range("D19").Value = "40.802.25"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("D21").Value = "'6.08"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "'71.26"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "'228.75"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("D25").Value = "'11.71"
$ws.Range("E25").Value = "  +14.86%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("E29").Value = "  -5.86%  "
$ws.Range("D30").Value = "'168.65"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("E31").Value = "  -8.75%  "
$ws.Range("D32").Value = "'20.14"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "'0.120"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "'5.62"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "'0.0739"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("D37").Value = "'4.53"
$ws.Range("E37").Value = "  -3.03%  "
$ws.Range("D38").Value = "'24.72"
$ws.Range("E38").Value = "  -6.35%  "
$ws.Range("D39").Value = "'3.95"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("E40").Value = "  +5.03%  "
$ws.Range("E41").Value = "  -4.86%  "
$ws.Range("D42").Value = "'5.44"
$ws.Range("E42").Value = "  -9.12%  "
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").Value = "'4.84"
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").Value = "'60.35"
$ws.Range("E45").Value = "  -12.27%  "
$ws.Range("E46").Value = "  -8.44%  "
$ws.Range("D47").Value = "'8.47"
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0993"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "'1.15"
$ws.Range("E51").Value = "  -3.43%  "
